$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "ID" sheet right before the "Model" sheet
# ------------------------------------------------------------------
$modelSheetRef = $wb.Worksheets.Item("Model")
$idSheet = $wb.Worksheets.Add($modelSheetRef)
$idSheet.Name = "ID"

# Worksheet references captured before the insert above can become stale,
# so re-fetch every sheet we still need to touch by name now.
$modelSheet = $wb.Worksheets.Item("Model")

# ------------------------------------------------------------------
# 2. Populate the new "ID" sheet
# ------------------------------------------------------------------
$idSheet.Range("A1").Value = "ID"
$idSheet.Range("B1").Value = "Rule"
$idSheet.Range("C1").Value = "Description"
$idSheet.Range("D1").Value = "Remarks"

$idSheet.Range("A2").Value = "TTYYNN"
$idSheet.Range("B2").Value = "TT (Model Type)"
$idSheet.Range("C2").Value = "모델 타입"
$idSheet.Range("D2").Value = "Model 시트 참조"

$idSheet.Range("B3").Value = "YY (Year)"
$idSheet.Range("C3").Value = "제품 출시 년도"

$idSheet.Range("B4").Value = "NN (Number)"
$idSheet.Range("C4").Value = "당해년도에 출시한 제품의 연변"

# Merge the ID rule cell spanning rows 2-4
$idSheet.Range("A2:A4").Merge()

# ------------------------------------------------------------------
# 3. Formatting for the new sheet
# ------------------------------------------------------------------

# Header row -> copy the look of the CompileOptions header (fill + border + centered)
$headerSrc = $wb.Worksheets.Item("CompileOptions").Range("A1")
$headerSrc.Copy()
$idSheet.Range("A1:D1").PasteSpecial(-4122)

# Body cells -> thin border around every populated cell
$bodyRange = $idSheet.Range("B2:D4")
$bodyRange.Borders.LineStyle = 1

# Merged ID cell -> thin border, centered horizontally and vertically
$idCell = $idSheet.Range("A2:A4")
$idCell.Borders.LineStyle = 1
$idCell.HorizontalAlignment = -4108
$idCell.VerticalAlignment = -4108

# Column widths (approximate target widths of 19.875 / 28.5 / 29.625 / 17.5 chars)
$idSheet.Columns.Item(1).ColumnWidth = 19.142857142857142
$idSheet.Columns.Item(2).ColumnWidth = 27.714285714285715
$idSheet.Columns.Item(3).ColumnWidth = 28.857142857142858
$idSheet.Columns.Item(4).ColumnWidth = 16.857142857142858

# Row heights
$idSheet.Rows.Item(1).RowHeight = 30
$idSheet.Rows.Item(2).RowHeight = 30
$idSheet.Rows.Item(3).RowHeight = 30
$idSheet.Rows.Item(4).RowHeight = 30

$idSheet.Range("B12").Select()

# ------------------------------------------------------------------
# 4. Update values on the CompileOptions sheet
# ------------------------------------------------------------------
$compile = $wb.Worksheets.Item("CompileOptions")
$compile.Range("A2").Value = 101900
$compile.Range("A3").Value = 101901
$compile.Range("A4").Value = 101902
$compile.Range("A5").Value = 101903

# ------------------------------------------------------------------
# 5. Update the selection on the Model sheet
# ------------------------------------------------------------------
$modelSheet.Range("B37").Select()

# ------------------------------------------------------------------
# 6. Restore the active sheet / selection to CompileOptions!A6
# ------------------------------------------------------------------
$compile.Activate()
$compile.Range("A6").Select()
